# Auto-generated Excel COM-interop script
# Applies a scheduled market-price data refresh across all leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), updating H/I/J/K/L/M/N price &
# profit columns to match freshly-pulled averages.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 823.8
$ws.Range("J4").Value = 1999
$ws.Range("L4").Value = 1999
$ws.Range("N4").Value = -2227
$ws.Range("H43").Value = 1744.75
$ws.Range("I43").Value = 1492.5
$ws.Range("J43").Value = 1997
$ws.Range("K43").Value = 1492.5
$ws.Range("L43").Value = 1997
$ws.Range("M43").Value = -1423.5
$ws.Range("N43").Value = -2135
$ws.Range("H132").Value = 2947.4167
$ws.Range("I132").Value = 2938.8965
$ws.Range("K132").Value = 8816.6895
$ws.Range("M132").Value = -6286.6895
$ws.Range("H135").Value = 3047.4075
$ws.Range("I135").Value = 2796.0476
$ws.Range("K135").Value = 25164.4284
$ws.Range("M135").Value = -22629.4284

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4154.2104
$ws.Range("I61").Value = 4159.2144
$ws.Range("J61").Value = 4140.2
$ws.Range("K61").Value = 4159.2144
$ws.Range("L61").Value = 4140.2
$ws.Range("M61").Value = -3947.2144
$ws.Range("N61").Value = -4564.2
$ws.Range("H74").Value = 3195.4546
$ws.Range("I74").Value = 2781.1428
$ws.Range("J74").Value = 3920.5
$ws.Range("K74").Value = 2781.1428
$ws.Range("L74").Value = 3920.5
$ws.Range("M74").Value = -1907.1428
$ws.Range("N74").Value = -5668.5
$ws.Range("H77").Value = 3195.4546
$ws.Range("I77").Value = 2781.1428
$ws.Range("J77").Value = 3920.5
$ws.Range("K77").Value = 13905.714
$ws.Range("L77").Value = 19602.5
$ws.Range("M77").Value = -9537.714
$ws.Range("N77").Value = -28338.5
$ws.Range("H97").Value = 2237.8518
$ws.Range("I97").Value = 1460.579
$ws.Range("K97").Value = 1460.579
$ws.Range("M97").Value = -964.579
$ws.Range("H101").Value = 54682
$ws.Range("J101").Value = 54682
$ws.Range("L101").Value = 54682
$ws.Range("N101").Value = -61172
$ws.Range("H102").Value = 15496.167
$ws.Range("I102").Value = 15496.167
$ws.Range("K102").Value = 15496.167
$ws.Range("M102").Value = -13874.167
$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040
$ws.Range("H132").Value = 1498.1515
$ws.Range("I132").Value = 1496.4517
$ws.Range("J132").Value = 1524.5
$ws.Range("K132").Value = 4489.355100000001
$ws.Range("L132").Value = 4573.5
$ws.Range("M132").Value = -1959.355100000001
$ws.Range("N132").Value = -9633.5
$ws.Range("H133").Value = 84884
$ws.Range("J133").Value = 84884
$ws.Range("L133").Value = 84884
$ws.Range("N133").Value = -89944
$ws.Range("H135").Value = 83119.39999999999
$ws.Range("J135").Value = 83119.39999999999
$ws.Range("L135").Value = 83119.39999999999
$ws.Range("N135").Value = -93259.39999999999
$ws.Range("H136").Value = 4154.2104
$ws.Range("I136").Value = 4159.2144
$ws.Range("J136").Value = 4140.2
$ws.Range("K136").Value = 12477.6432
$ws.Range("L136").Value = 12420.6
$ws.Range("M136").Value = -9927.643199999999
$ws.Range("N136").Value = -17520.6

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1561.8
$ws.Range("I99").Value = 705.25
$ws.Range("J99").Value = 4988
$ws.Range("K99").Value = 705.25
$ws.Range("L99").Value = 4988
$ws.Range("M99").Value = 792.75
$ws.Range("N99").Value = -7984
$ws.Range("H105").Value = 3167.3438
$ws.Range("I105").Value = 2871.7407
$ws.Range("J105").Value = 4763.6
$ws.Range("K105").Value = 2871.7407
$ws.Range("L105").Value = 4763.6
$ws.Range("M105").Value = -1124.7407
$ws.Range("N105").Value = -8257.6

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 946.2857
$ws.Range("J22").Value = 1077.25
$ws.Range("L22").Value = 1077.25
$ws.Range("N22").Value = -1777.25
$ws.Range("H31").Value = 4165.45
$ws.Range("I31").Value = 3543.3076
$ws.Range("J31").Value = 5320.857
$ws.Range("K31").Value = 3543.3076
$ws.Range("L31").Value = 5320.857
$ws.Range("M31").Value = -3248.3076
$ws.Range("N31").Value = -5910.857
$ws.Range("H34").Value = 4165.45
$ws.Range("I34").Value = 3543.3076
$ws.Range("J34").Value = 5320.857
$ws.Range("K34").Value = 3543.3076
$ws.Range("L34").Value = 5320.857
$ws.Range("M34").Value = -3341.3076
$ws.Range("N34").Value = -5724.857
$ws.Range("H58").Value = 1907
$ws.Range("I58").Value = 884.125
$ws.Range("J58").Value = 5180.2
$ws.Range("K58").Value = 884.125
$ws.Range("L58").Value = 5180.2
$ws.Range("M58").Value = -681.125
$ws.Range("N58").Value = -5586.2
$ws.Range("H99").Value = 1612.1666
$ws.Range("I99").Value = 1671.7778
$ws.Range("J99").Value = 1433.3334
$ws.Range("K99").Value = 1671.7778
$ws.Range("L99").Value = 1433.3334
$ws.Range("M99").Value = -173.7778000000001
$ws.Range("N99").Value = -4429.3334
$ws.Range("H105").Value = 3147.8333
$ws.Range("I105").Value = 2642
$ws.Range("K105").Value = 2642
$ws.Range("M105").Value = -895
$ws.Range("H119").Value = 69760.5
$ws.Range("J119").Value = 69760.5
$ws.Range("L119").Value = 69760.5
$ws.Range("N119").Value = -79436.5
$ws.Range("H126").Value = 1612.1666
$ws.Range("I126").Value = 1671.7778
$ws.Range("J126").Value = 1433.3334
$ws.Range("K126").Value = 5015.3334
$ws.Range("L126").Value = 4300.0002
$ws.Range("M126").Value = -2545.3334
$ws.Range("N126").Value = -9240.0002
$ws.Range("H132").Value = 2536.8125
$ws.Range("I132").Value = 2235
$ws.Range("K132").Value = 6705
$ws.Range("M132").Value = -4175
$ws.Range("H134").Value = 1099.75
$ws.Range("I134").Value = 999.7646999999999
$ws.Range("J134").Value = 1666.3334
$ws.Range("K134").Value = 2999.2941
$ws.Range("L134").Value = 4999.0002
$ws.Range("M134").Value = -464.2941000000001
$ws.Range("N134").Value = -10069.0002
$ws.Range("H136").Value = 1907
$ws.Range("I136").Value = 884.125
$ws.Range("J136").Value = 5180.2
$ws.Range("K136").Value = 2652.375
$ws.Range("L136").Value = 15540.6
$ws.Range("M136").Value = -102.375
$ws.Range("N136").Value = -20640.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 728.64703
$ws.Range("I5").Value = 499.18182
$ws.Range("K5").Value = 1497.54546
$ws.Range("M5").Value = -1385.54546
$ws.Range("H131").Value = 18519734
$ws.Range("I131").Value = 7937448
$ws.Range("J131").Value = 33334934
$ws.Range("K131").Value = 23812344
$ws.Range("L131").Value = 100004802
$ws.Range("M131").Value = -23807304
$ws.Range("N131").Value = -100014882
$ws.Range("H135").Value = 728.64703
$ws.Range("I135").Value = 499.18182
$ws.Range("K135").Value = 4492.63638
$ws.Range("M135").Value = -1957.63638

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 62.333332
$ws.Range("I2").Value = 62.333332
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 62.333332
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 50.666668
$ws.Range("N2").ClearContents()
$ws.Range("H101").Value = 32916.75
$ws.Range("J101").Value = 32916.75
$ws.Range("L101").Value = 32916.75
$ws.Range("N101").Value = -39406.75
$ws.Range("H102").Value = 10454.969
$ws.Range("I102").Value = 12192.167
$ws.Range("K102").Value = 12192.167
$ws.Range("M102").Value = -10570.167
$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 240000
$ws.Range("N134").Value = -245070
$ws.Range("H136").Value = 28517.643
$ws.Range("J136").Value = 28517.643
$ws.Range("L136").Value = 85552.929
$ws.Range("N136").Value = -90652.929

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2254.1333
$ws.Range("I40").Value = 1928.3
$ws.Range("J40").Value = 2905.8
$ws.Range("K40").Value = 1928.3
$ws.Range("L40").Value = 2905.8
$ws.Range("M40").Value = -1792.3
$ws.Range("N40").Value = -3177.8
$ws.Range("H46").Value = 1797.2307
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 2062.6667
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 2062.6667
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -2438.6667
$ws.Range("H55").Value = 439.4737
$ws.Range("J55").Value = 590.44446
$ws.Range("L55").Value = 590.44446
$ws.Range("N55").Value = -936.44446
$ws.Range("H132").Value = 2843.0322
$ws.Range("I132").Value = 2268.3076
$ws.Range("K132").Value = 6804.9228
$ws.Range("M132").Value = -4274.9228

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 58020
$ws.Range("J46").Value = 58020
$ws.Range("L46").Value = 58020
$ws.Range("N46").Value = -58482
$ws.Range("H62").Value = 8984
$ws.Range("J62").Value = 8988.5
$ws.Range("L62").Value = 8988.5
$ws.Range("N62").Value = -10236.5
$ws.Range("H65").Value = 8984
$ws.Range("J65").Value = 8988.5
$ws.Range("L65").Value = 44942.5
$ws.Range("N65").Value = -51182.5
$ws.Range("H132").Value = 1888.8684
$ws.Range("J132").Value = 2997.8572
$ws.Range("L132").Value = 8993.571599999999
$ws.Range("N132").Value = -14053.5716
$ws.Range("H134").Value = 58020
$ws.Range("J134").Value = 58020
$ws.Range("L134").Value = 174060
$ws.Range("N134").Value = -179130
